$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

# --- Sheet: Therbligs(L) -- extend Gantt rows 10:17 (mirrors rows 2:9) ---
$ws2 = $wb.Worksheets.Item("Therbligs(L)")
$ws2.Range("A10").Value = "R"
$ws2.Range("A2").Copy() | Out-Null
$ws2.Range("A10").PasteSpecial($xlPasteFormats)
$ws2.Range("B10").Value = "AGENT"
$ws2.Range("B2").Copy() | Out-Null
$ws2.Range("B10").PasteSpecial($xlPasteFormats)
$ws2.Range("C10").Value = "Pillar"
$ws2.Range("C2").Copy() | Out-Null
$ws2.Range("C10").PasteSpecial($xlPasteFormats)
$ws2.Range("D10").Value = "B"
$ws2.Range("D2").Copy() | Out-Null
$ws2.Range("D10").PasteSpecial($xlPasteFormats)
$ws2.Range("A11").Value = "G"
$ws2.Range("A3").Copy() | Out-Null
$ws2.Range("A11").PasteSpecial($xlPasteFormats)
$ws2.Range("E11").Value = "Pillar"
$ws2.Range("E3").Copy() | Out-Null
$ws2.Range("E11").PasteSpecial($xlPasteFormats)
$ws2.Range("A12").Value = "P"
$ws2.Range("A4").Copy() | Out-Null
$ws2.Range("A12").PasteSpecial($xlPasteFormats)
$ws2.Range("E12").Value = "Pillar"
$ws2.Range("E4").Copy() | Out-Null
$ws2.Range("E12").PasteSpecial($xlPasteFormats)
$ws2.Range("F12").Value = "BottomPlate"
$ws2.Range("F4").Copy() | Out-Null
$ws2.Range("F12").PasteSpecial($xlPasteFormats)
$ws2.Range("A13").Value = "RL"
$ws2.Range("A5").Copy() | Out-Null
$ws2.Range("A13").PasteSpecial($xlPasteFormats)
$ws2.Range("A14").Value = "R"
$ws2.Range("A6").Copy() | Out-Null
$ws2.Range("A14").PasteSpecial($xlPasteFormats)
$ws2.Range("B14").Value = "Pillar"
$ws2.Range("B6").Copy() | Out-Null
$ws2.Range("B14").PasteSpecial($xlPasteFormats)
$ws2.Range("C14").Value = "BottomPlate"
$ws2.Range("C6").Copy() | Out-Null
$ws2.Range("C14").PasteSpecial($xlPasteFormats)
$ws2.Range("D14").Value = "B"
$ws2.Range("D6").Copy() | Out-Null
$ws2.Range("D14").PasteSpecial($xlPasteFormats)
$ws2.Range("A15").Value = "G"
$ws2.Range("A7").Copy() | Out-Null
$ws2.Range("A15").PasteSpecial($xlPasteFormats)
$ws2.Range("E15").Value = "Pillar"
$ws2.Range("E7").Copy() | Out-Null
$ws2.Range("E15").PasteSpecial($xlPasteFormats)
$ws2.Range("A16").Value = "P"
$ws2.Range("A8").Copy() | Out-Null
$ws2.Range("A16").PasteSpecial($xlPasteFormats)
$ws2.Range("E16").Value = "Pillar"
$ws2.Range("E8").Copy() | Out-Null
$ws2.Range("E16").PasteSpecial($xlPasteFormats)
$ws2.Range("F16").Value = "BottomPlate"
$ws2.Range("F8").Copy() | Out-Null
$ws2.Range("F16").PasteSpecial($xlPasteFormats)
$ws2.Range("A17").Value = "RL"
$ws2.Range("A9").Copy() | Out-Null
$ws2.Range("A17").PasteSpecial($xlPasteFormats)
$ws2.Activate()
$ws2.Range("A10:F17").Select()

# --- Sheet: Therbligs(R) -- extend Gantt rows 8:13 (mirrors rows 2:7) ---
$ws3 = $wb.Worksheets.Item("Therbligs(R)")
$ws3.Range("A8").Value = "R"
$ws3.Range("A2").Copy() | Out-Null
$ws3.Range("A8").PasteSpecial($xlPasteFormats)
$ws3.Range("B8").Value = "AGENT"
$ws3.Range("B2").Copy() | Out-Null
$ws3.Range("B8").PasteSpecial($xlPasteFormats)
$ws3.Range("C8").Value = "BottomPlate"
$ws3.Range("C2").Copy() | Out-Null
$ws3.Range("C8").PasteSpecial($xlPasteFormats)
$ws3.Range("D8").Value = "B"
$ws3.Range("D2").Copy() | Out-Null
$ws3.Range("D8").PasteSpecial($xlPasteFormats)
$ws3.Range("A9").Value = "G"
$ws3.Range("A3").Copy() | Out-Null
$ws3.Range("A9").PasteSpecial($xlPasteFormats)
$ws3.Range("E9").Value = "Pillar"
$ws3.Range("E3").Copy() | Out-Null
$ws3.Range("E9").PasteSpecial($xlPasteFormats)
$ws3.Range("A10").Value = "M"
$ws3.Range("A4").Copy() | Out-Null
$ws3.Range("A10").PasteSpecial($xlPasteFormats)
$ws3.Range("B10").Value = "BottomPlate"
$ws3.Range("B4").Copy() | Out-Null
$ws3.Range("B10").PasteSpecial($xlPasteFormats)
$ws3.Range("C10").Value = "BottomPlate"
$ws3.Range("C4").Copy() | Out-Null
$ws3.Range("C10").PasteSpecial($xlPasteFormats)
$ws3.Range("D10").Value = "B"
$ws3.Range("D4").Copy() | Out-Null
$ws3.Range("D10").PasteSpecial($xlPasteFormats)
$ws3.Range("A11").Value = "P"
$ws3.Range("A5").Copy() | Out-Null
$ws3.Range("A11").PasteSpecial($xlPasteFormats)
$ws3.Range("E11").Value = "Pillar"
$ws3.Range("E5").Copy() | Out-Null
$ws3.Range("E11").PasteSpecial($xlPasteFormats)
$ws3.Range("F11").Value = "BottomPlate"
$ws3.Range("F5").Copy() | Out-Null
$ws3.Range("F11").PasteSpecial($xlPasteFormats)
$ws3.Range("A12").Value = "A"
$ws3.Range("A6").Copy() | Out-Null
$ws3.Range("A12").PasteSpecial($xlPasteFormats)
$ws3.Range("E12").Value = "Pillar"
$ws3.Range("E6").Copy() | Out-Null
$ws3.Range("E12").PasteSpecial($xlPasteFormats)
$ws3.Range("F12").Value = "BottomPlate"
$ws3.Range("F6").Copy() | Out-Null
$ws3.Range("F12").PasteSpecial($xlPasteFormats)
$ws3.Range("A13").Value = "RL"
$ws3.Range("A7").Copy() | Out-Null
$ws3.Range("A13").PasteSpecial($xlPasteFormats)
$ws3.Range("E13").Value = "Pillar"
$ws3.Range("E7").Copy() | Out-Null
$ws3.Range("E13").PasteSpecial($xlPasteFormats)
$ws3.Activate()
$ws3.Range("A8:F13").Select()

# --- Sheet: OHT Relation -- extend rk-evolution matrix to A1:H8 ---
$ws5 = $wb.Worksheets.Item("OHT Relation")
$ws5.Range("F1").Value = 4
$ws5.Range("G1").Value = 5
$ws5.Range("H1").Value = 6
$ws5.Range("F2").Value = 0
$ws5.Range("G2").Value = 0
$ws5.Range("H2").Value = 0
$ws5.Range("D3").Value = 1
$ws5.Range("E3").Value = 0
$ws5.Range("F3").Value = 0
$ws5.Range("G3").Value = 0
$ws5.Range("H3").Value = 0
$ws5.Range("C4").Value = -1
$ws5.Range("E4").Value = 0
$ws5.Range("F4").Value = 0
$ws5.Range("G4").Value = 0
$ws5.Range("H4").Value = 1
$ws5.Range("C5").Value = 0
$ws5.Range("D5").Value = 0
$ws5.Range("F5").Value = 1
$ws5.Range("G5").Value = 1
$ws5.Range("H5").Value = 0
$ws5.Range("A6").Value = 4
$ws5.Range("B6").Value = 0
$ws5.Range("C6").Value = 0
$ws5.Range("D6").Value = 0
$ws5.Range("E6").Value = -1
$ws5.Range("F6").Value = 0
$ws5.Range("G6").Value = 0
$ws5.Range("H6").Value = 1
$ws5.Range("A7").Value = 5
$ws5.Range("B7").Value = 0
$ws5.Range("C7").Value = 0
$ws5.Range("D7").Value = 0
$ws5.Range("E7").Value = -1
$ws5.Range("F7").Value = 0
$ws5.Range("G7").Value = 0
$ws5.Range("H7").Value = 1
$ws5.Range("A8").Value = 6
$ws5.Range("B8").Value = 0
$ws5.Range("C8").Value = 0
$ws5.Range("D8").Value = -1
$ws5.Range("E8").Value = 0
$ws5.Range("F8").Value = -1
$ws5.Range("G8").Value = -1
$ws5.Range("H8").Value = 0
$ws5.Range("H8").Select()

Write-Host "edit applied"
